# Tabu search partially repaired.
# Update allocation values (Barbera=B, Chardonnay=C, Nebbiolo=D) for three
# monthly sheets, moving allocations between grape varieties for a few plots.

$wb = $excel.ActiveWorkbook

# --- Miesiac 3 ---
$ws3 = $wb.Worksheets.Item("Miesiac 3")
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 138
$ws3.Range("C3").Value = 0
$ws3.Range("D3").Value = 115
$ws3.Range("C4").Value = 265
$ws3.Range("D4").Value = 0

# --- Miesiac 7 ---
$ws7 = $wb.Worksheets.Item("Miesiac 7")
$ws7.Range("B2").Value = 0
$ws7.Range("C2").Value = 124
$ws7.Range("B3").Value = 0
$ws7.Range("C3").Value = 100
$ws7.Range("C4").Value = 166
$ws7.Range("D4").Value = 0

# --- Miesiac 11 ---
$ws11 = $wb.Worksheets.Item("Miesiac 11")
$ws11.Range("B2").Value = 178
$ws11.Range("D2").Value = 0
$ws11.Range("B3").Value = 190
$ws11.Range("C3").Value = 0
$ws11.Range("B4").Value = 213
